$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2024-04-24 Wednesday" "2024-04-25 Thursday"

Replace-Text "79×80=" "70×97="
Replace-Text "52×77=" "38×61="
Replace-Text "75×51=" "73×36="
Replace-Text "49×80=" "74×70="
Replace-Text "70×53=" "58×51="

Replace-Text "46×63=" "11×94="
Replace-Text "33×20=" "97×27="
Replace-Text "37×50=" "38×46="
Replace-Text "33×71=" "23×85="
Replace-Text "21×17=" "62×13="

Replace-Text "14×82=" "47×80="
Replace-Text "74×11=" "57×81="
Replace-Text "37×27=" "52×90="
Replace-Text "50×18=" "47×87="
Replace-Text "74×94=" "38×14="

Replace-Text "33×59=" "74×58="
Replace-Text "12×92=" "29×97="
Replace-Text "47×27=" "16×28="
Replace-Text "18×92=" "29×76="
Replace-Text "32×28=" "67×36="

Replace-Text "61×75=" "18×23="
Replace-Text "81×41=" "98×43="
Replace-Text "82×57=" "61×38="
Replace-Text "79×14=" "95×90="
Replace-Text "71×42=" "28×75="
